$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''26.289.93'
$ws.Range("E2").Value = '''  +0.21%  '
$ws.Range("D3").Value = '''1.680.27'
$ws.Range("E3").Value = '''  +0.47%  '
$ws.Range("D4").Value = '''1.007'
$ws.Range("E4").Value = '''  +0.18%  '
$ws.Range("D5").Value = '''218.38'
$ws.Range("E5").Value = '''  +0.20%  '
$ws.Range("D6").Value = '''0.5273'
$ws.Range("E6").Value = '''  +2.58%  '
$ws.Range("E7").Value = '''  +0.17%  '
$ws.Range("D8").Value = '''0.2700'
$ws.Range("E8").Value = '''  +1.42%  '
$ws.Range("E9").Value = '''  +0.70%  '
$ws.Range("D10").Value = '''22.02'
$ws.Range("E10").Value = '''  +1.97%  '
$ws.Range("D11").Value = '''0.07510'
$ws.Range("E11").Value = '''  +1.76%  '
$ws.Range("D12").Value = '''4.544'
$ws.Range("E12").Value = '''  -0.15%  '
$ws.Range("D13").Value = '''1.676.66'
$ws.Range("E13").Value = '''  +0.28%  '
$ws.Range("D14").Value = '''0.5804'
$ws.Range("E14").Value = '''  -0.59%  '
$ws.Range("D15").Value = '''0.000008477'
$ws.Range("E15").Value = '''  -2.17%  '
$ws.Range("D16").Value = '''64.26'
$ws.Range("E16").Value = '''  -0.71%  '
$ws.Range("D17").Value = '''26.308.36'
$ws.Range("E17").Value = '''  +0.05%  '
$ws.Range("E18").Value = '''  -0.93%  '
$ws.Range("E19").Value = '''  +0.11%  '
$ws.Range("D20").Value = '''10.87'
$ws.Range("E20").Value = '''  +0.00%  '
$ws.Range("D21").Value = '''189.40'
$ws.Range("E21").Value = '''  -0.10%  '
$ws.Range("D22").Value = '''6.205'
$ws.Range("D23").Value = '''1.008'
$ws.Range("E23").Value = '''  +0.11%  '
$ws.Range("D24").Value = '''144.90'
$ws.Range("E24").Value = '''  +0.32%  '
$ws.Range("D25").Value = '''7.718'
$ws.Range("E25").Value = '''  +1.00%  '
$ws.Range("D26").Value = '''0.1238'
$ws.Range("E26").Value = '''  +4.38%  '
$ws.Range("E27").Value = '''  +0.92%  '
$ws.Range("D28").Value = '''0.06587'
$ws.Range("E28").Value = '''  +10.23%  '
$ws.Range("D29").Value = '''1.361'
$ws.Range("E29").Value = '''  +6.08%  '
$ws.Range("D30").Value = '''1.326'
$ws.Range("E30").Value = '''  -0.02%  '
$ws.Range("D31").Value = '''3.581'
$ws.Range("E31").Value = '''  +1.37%  '
$ws.Range("E32").Value = '''  +1.05%  '
$ws.Range("D33").Value = '''1.660'
$ws.Range("E33").Value = '''  +1.15%  '
$ws.Range("E34").Value = '''  +0.89%  '
$ws.Range("D35").Value = '''0.6188'
$ws.Range("E35").Value = '''  +2.68%  '
$ws.Range("D36").Value = '''2.397'
$ws.Range("E36").Value = '''  +0.97%  '
$ws.Range("D37").Value = '''2.700'
$ws.Range("E37").Value = '''  +1.99%  '
$ws.Range("D38").Value = '''6.372'
$ws.Range("E38").Value = '''  +4.77%  '
$ws.Range("E39").Value = '''  -0.02%  '
$ws.Range("D40").Value = '''1.107.34'
$ws.Range("E40").Value = '''  +2.40%  '
$ws.Range("D41").Value = '''0.8761'
$ws.Range("E41").Value = '''  +0.62%  '
$ws.Range("D42").Value = '''1.014'
$ws.Range("E42").Value = '''  +0.35%  '
$ws.Range("D43").Value = '''100.51'
$ws.Range("E43").Value = '''  +0.30%  '
$ws.Range("D44").Value = '''1.828.04'
$ws.Range("D45").Value = '''0.00000000111'
$ws.Range("E45").Value = '''  -2.56%  '
$ws.Range("D46").Value = '''56.83'
$ws.Range("E46").Value = '''  +1.14%  '
$ws.Range("B47").Value = '''Frax'
$ws.Range("C47").Value = '''https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D47").Value = '''1.008'
$ws.Range("E47").Value = '''  -0.41%  '
$ws.Range("B48").Value = '''EnergySwap'
$ws.Range("C48").Value = '''https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '''8.142'
$ws.Range("E48").Value = '''  +0.91%  '
$ws.Range("D49").Value = '''0.05267'
$ws.Range("E49").Value = '''  +1.05%  '
$ws.Range("D50").Value = '''0.4302'
$ws.Range("E50").Value = '''  +0.14%  '
$ws.Range("E51").Value = '''  +2.45%  '
